# vydani verze 3.7.0, plne funkcni ip manager, nova vizualizace menu
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ip_address_list": keep only row 1 (edited), drop rows 2-6
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ip_address_list")
$ws1.Activate() | Out-Null

$ws1.Range("A2:E6").Delete() | Out-Null

$ws1.Range("A1").Value = "fff"
$ws1.Range("B1").Value = "192.168.10.240"
$ws1.Range("D1").ClearContents() | Out-Null
$ws1.Range("E1").Value = $false

$ws1.Range("G11").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "ip_adress_fav_list": wipe all rows, sheet becomes empty
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")
$ws2.Activate() | Out-Null

$ws2.Range("A1:E3").Delete() | Out-Null

$ws2.Range("A1:E3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "disk_list": add multi-line note in F3, update selection, becomes
# the active/selected tab
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("disk_list")
$ws3.Activate() | Out-Null

$ws3.Range("F3").Value = "a
a
a
a
aa"

$ws3.Range("C16").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Settings": update B1 value 1 -> 4
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Settings")
$ws4.Range("B1").Value = 4

# ---------------------------------------------------------------------------
# Re-activate "disk_list" last so it ends up as the selected/visible tab
# (tabSelected moves off "Settings" and onto "disk_list"; activeTab -> 2)
# ---------------------------------------------------------------------------
$ws3.Activate() | Out-Null

# Best-effort: restore/refresh the workbook window geometry recorded in the
# file's bookViews (engine may not persist these via COM, kept for parity).
$win = $excel.ActiveWindow
$win.Left = 14940
$win.Top = 2640
$win.Width = 23010
$win.Height = 13650
